$d = $word.ActiveDocument

# 1. Remove the word "but " from the "Keeping up to speed..." paragraph.
$d.Content.Find.Execute(
    "provide great efficiencies, but they require",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "provide great efficiencies, they require", 2)

# 2. Merge the standalone " " run with the following "of the software" run
#    (removes a redundant run split without changing the visible text).
$r = $d.Content
$r.Find.Execute("of the software", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($r.Start, $r.End)
$target.Find.Execute(
    "of the software",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "of the software", 2)

# 3. Remove the "_GoBack" bookmark left over at the end of the document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
